# "final files for 6/7"
# Adds a new "Predictive modeling" category (3 rows) to the Analysis
# packages list, between "Power Analysis"/"PUMP" (row 28) and
# "Propensity Score Matching" (old row 29), and un-merges / renumbers the
# old "Structural Equation Modeling" / "SEMinR packagez" pair so the
# second row (SEMinR) is no longer merged with the first (sem package).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert three new blank rows right after row 28 ("PUMP ...") and
#    before the old row 29 ("Propensity Score Matching"). Excel shifts
#    everything below down automatically (row refs, merged cells and
#    hyperlinks all re-point to the new rows).
# ------------------------------------------------------------------
$ws.Rows("29:31").Insert()

# ------------------------------------------------------------------
# 2) Column A: new category label "Predictive modeling", merged over
#    the 3 new rows (same look as the other category cells, e.g. the
#    "Power Analysis" / "Network Analysis" blocks just above it).
# ------------------------------------------------------------------
$ws.Range("A26").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Predictive modeling"
$ws.Range("A29:A31").Merge()

# ------------------------------------------------------------------
# 3) Column B: the three package/resource entries for the new
#    category, formatted like the other hyperlink-style entries.
#    Values are entered in the same order the original author typed
#    them (Tidymodels, then caret, then the "lots of packages" row)
#    so the shared-string table comes out in the same order.
# ------------------------------------------------------------------
$ws.Range("B26").Copy()
$ws.Range("B29:B31").PasteSpecial(-4122)

$ws.Range("B31").Value = "Tidymodels (part of Tidyverse)"
$ws.Range("B29").Value = "caret"
$ws.Range("B30").Value = "lots of packages for individual maching learning algorithms (e.g. randomForest)"

# Hyperlinks - added in the same order the original author used
# (caret, then Tidymodels, then the "lots of packages" row).
$ws.Hyperlinks.Add($ws.Range("B29"), "https://cran.r-project.org/web/packages/caret/index.html")
$ws.Hyperlinks.Add($ws.Range("B31"), "https://cran.r-project.org/web/packages/tidymodels/index.html")
$ws.Hyperlinks.Add($ws.Range("B30"), "https://cran.r-project.org/web/packages/randomForest/index.html")

# Re-apply the hyperlink-style formatting (Hyperlinks.Add resets the
# cell style to its own default blue/underline style object).
$ws.Range("B26").Copy()
$ws.Range("B29:B31").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4) Un-merge the old "Structural Equation Modeling" / "SEMinR
#    packagez" pair (now at rows 45:46 after the insert above) so
#    each row keeps its own plain-category style instead of sharing
#    the old merged block.
# ------------------------------------------------------------------
$ws.Range("A45:A46").UnMerge()
$ws.Range("A44").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A46").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 5) Update the view state to match where the author left off:
#    scrolled down a bit further and the new category selected.
# ------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A17"))
$ws.Range("A29:A31").Select()
